$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions) -- F column ("想去人数") numeric refresh
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 7875
$wsExpo.Range("F5").Value = 7875
$wsExpo.Range("F9").Value = 8642
$wsExpo.Range("F10").Value = 8642
$wsExpo.Range("F13").Value = 90
$wsExpo.Range("F14").Value = 5777
$wsExpo.Range("F16").Value = 2768
$wsExpo.Range("F17").Value = 1202
$wsExpo.Range("F20").Value = 43
$wsExpo.Range("F22").Value = 80
$wsExpo.Range("F23").Value = 3928
$wsExpo.Range("F24").Value = 81
$wsExpo.Range("F25").Value = 64
$wsExpo.Range("F28").Value = 182
$wsExpo.Range("F30").Value = 5514
$wsExpo.Range("F31").Value = 9
$wsExpo.Range("F34").Value = 393
$wsExpo.Range("F35").Value = 162
$wsExpo.Range("F36").Value = 397
$wsExpo.Range("F37").Value = 2593
$wsExpo.Range("F38").Value = 1528
$wsExpo.Range("F41").Value = 4968
$wsExpo.Range("F45").Value = 3598
$wsExpo.Range("F46").Value = 9
$wsExpo.Range("F50").Value = 473
$wsExpo.Range("F51").Value = 21

# ---------------------------------------------------------------------
# Sheet "演出" (Performances) -- insert a new event row (2024-05-11,
# 赵鹏 concert) above the existing row 8, shifting rows 8-10 down to
# 9-11, then fix up the running index + the F-column ("想去人数") values.
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Rows.Item(8).Insert()

$wsShow.Range("A8").Value = 7
# Force text (not an auto-converted date serial) to match the rest of
# column B, which stores these as literal "yyyy-mm-dd" strings.
$wsShow.Range("B8").NumberFormat = "@"
$wsShow.Range("B8").Value = "2024-05-11"
$wsShow.Range("C8").Value = "北京·赵鹏“行吟2024·拥抱”——巡演十周年纪念演唱会"
$wsShow.Range("D8").Value = "西直门外大街135号（北京展览馆内） 北京展览馆剧场"
$wsShow.Range("E8").Value = "2024.05.11 19:30-05.11 21:00"
$wsShow.Range("F8").Value = 0
$wsShow.Range("G8").Value = 299
$wsShow.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=82859"
$wsShow.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202403/t9FYkHQt1710397403533.jpeg"

# The insert pushed the old rows 8-10 down to 9-11, carrying their old
# running-index values (7, 8, 9) along for the ride -- bump each back
# onto the correct sequence (8, 9, 10).
$wsShow.Range("A9").Value = 8
$wsShow.Range("A10").Value = 9
$wsShow.Range("A11").Value = 10

# Unrelated "想去人数" refreshes elsewhere on the same sheet (rows 3 and 5
# are untouched by the insert; row 10's count -- the "剧院魅影" show, now
# shifted down from row 9 -- also ticked up by one).
$wsShow.Range("F3").Value = 153
$wsShow.Range("F5").Value = 73
$wsShow.Range("F10").Value = 131

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types) -- F column ("想去人数") numeric refresh
# (this combined view is not re-paginated by the new 演出 row, only the
# same underlying counters move)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 7875
$wsAll.Range("F6").Value = 7875
$wsAll.Range("F9").Value = 8642
$wsAll.Range("F10").Value = 8642
$wsAll.Range("F12").Value = 90
$wsAll.Range("F13").Value = 5777
$wsAll.Range("F15").Value = 2768
$wsAll.Range("F16").Value = 1202
$wsAll.Range("F19").Value = 43
$wsAll.Range("F20").Value = 153
$wsAll.Range("F22").Value = 80
$wsAll.Range("F23").Value = 3928
$wsAll.Range("F24").Value = 81
$wsAll.Range("F25").Value = 64
$wsAll.Range("F28").Value = 182
$wsAll.Range("F30").Value = 5515
$wsAll.Range("F31").Value = 9
$wsAll.Range("F33").Value = 393
$wsAll.Range("F34").Value = 162
$wsAll.Range("F35").Value = 397
$wsAll.Range("F36").Value = 73
$wsAll.Range("F37").Value = 2593
$wsAll.Range("F38").Value = 1528
$wsAll.Range("F42").Value = 4968
$wsAll.Range("F46").Value = 3598
$wsAll.Range("F48").Value = 473
$wsAll.Range("F49").Value = 21
$wsAll.Range("F50").Value = 131